# CIERRE 24 JUN 2022
# Advance the payroll receipt sheet from "SEMANA 24" (13-19 Jun 2022) to
# "SEMANA 25" (20-26 Jun 2022): update the week-label text, the two
# "extra"/discount amounts for the new week, and move the on-screen
# scroll/selection down to where the next week's entry rows are.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# 1) Week label (B9). H9 (=B9), B27 (=B9), H27 (=B27) and B43 (=H27) are
#    formulas that recompute automatically from this single source cell.
$ws.Range("B9").Value = "SEMANA   25  DEL    20      Al   26   DE   JUNIO          2022"

# 2) Updated figures for the new week.
$ws.Range("K21").Value = 1120   # K24 = SUM(K21:K23) recalculates automatically
$ws.Range("E40").Value = 1250   # E41 = SUM(E38:E40) recalculates automatically

# 3) Move the visible window / selection to the next block of entry rows.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H39:I40").Select() | Out-Null
